# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" worker table (rows 16-28) contained obsolete
# records for periods 1605/1606 belonging to 4 workers (94483278,
# 1047427860, 1143338624, 80812126), followed by 5 up-to-date records for
# period 1612 (1116435458, 1118288813, 22790491, 14899808, 1047446376).
# The obsolete rows must be removed so only the 5 current-period workers
# remain, and the summary totals (total "Valor Mora", worker count and
# period count) must be refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 (the last obsolete record) carries the special "bottom of table"
# border formatting. Row 24 (ANGEL DE JESUS MESTRA ZULETA / period 1612)
# is the record that will become the new last row of the table once the
# obsolete rows below it are deleted, so give it that formatting first.
$ws.Range("B28:J28").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Delete the obsolete rows, working from the bottom up so the row numbers
# of the rows still waiting to be deleted don't shift underneath us:
#   25-28: 1143338624 PEDRO LUIS PULIDO ANAYA (periods 1606 & 1605) and
#          80812126 HENIEK YITZAK HERRERA RAMIREZ (periods 1606 & 1605)
#   16-19: 94483278 EIDER PIMENTEL CALDON (periods 1606 & 1605) and
#          1047427860 YINETH PAOLA BALLESTAS FERIA (periods 1606 & 1605)
$ws.Rows("25:28").Delete()
$ws.Rows("16:19").Delete()

# Refresh the summary figures at the top of the statement to match the
# remaining 5 workers / 1 period of data.
$ws.Range("E11").Value = 226320
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 1
